$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The handback transform failed for the c4c4bfd1 file (row 3 everywhere) -
# update the "Status" column on the Overview sheet and on each language sheet.
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file name mismatch in the "Error Detail" column
# (column K) for row 3 on both language sheets.
$zhcn.Range("K3").Value = "Handback file name: t5fsvlz0.lbe is different with handoff file name: c4c4bfd1-2999-4c79-8291-89c1bd78cde6.04dfe7bc2290248f01550edbbaaa6698596755b6.zh-cn."
$dede.Range("K3").Value = "Handback file name: t5fsvlz0.lbe is different with handoff file name: c4c4bfd1-2999-4c79-8291-89c1bd78cde6.04dfe7bc2290248f01550edbbaaa6698596755b6.de-de."
